$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> corrected numeric value for the naive component
# forecaster bug fix (rows 24-53, columns B:K).
$values = @{
    "B24" = -0.4
    "C24" = -0.7
    "D24" = 0.6
    "E24" = 0.6
    "F24" = -0.3
    "G24" = 0.7
    "H24" = -0.6
    "I24" = -0.3
    "J24" = -0.3
    "K24" = 2
    "B25" = -0.6699999999999999
    "C25" = 0.63
    "D25" = 0.63
    "E25" = -0.27
    "F25" = 0.73
    "G25" = -0.57
    "H25" = -0.27
    "I25" = -0.27
    "J25" = 2.03
    "K25" = -1.37
    "B26" = 0.72
    "C26" = 0.72
    "D26" = -0.18
    "E26" = 0.82
    "F26" = -0.48
    "G26" = -0.18
    "H26" = -0.18
    "I26" = 2.12
    "J26" = -1.28
    "K26" = -1.28
    "B27" = 0.6899999999999999
    "C27" = -0.21
    "D27" = 0.7899999999999999
    "E27" = -0.51
    "F27" = -0.21
    "G27" = -0.21
    "H27" = 2.09
    "I27" = -1.31
    "J27" = -1.31
    "K27" = -0.01000000000000001
    "B28" = -0.29
    "C28" = 0.71
    "D28" = -0.59
    "E28" = -0.29
    "F28" = -0.29
    "G28" = 2.01
    "H28" = -1.39
    "I28" = -1.39
    "J28" = -0.08999999999999997
    "K28" = 0.11
    "B29" = 0.6499999999999999
    "C29" = -0.65
    "D29" = -0.35
    "E29" = -0.35
    "F29" = 1.95
    "G29" = -1.45
    "H29" = -1.45
    "I29" = -0.15
    "J29" = 0.05
    "K29" = -0.15
    "B30" = -0.6699999999999999
    "C30" = -0.37
    "D30" = -0.37
    "E30" = 1.93
    "F30" = -1.47
    "G30" = -1.47
    "H30" = -0.17
    "I30" = 0.03000000000000001
    "J30" = -0.17
    "K30" = -1.27
    "B31" = -0.36
    "C31" = -0.36
    "D31" = 1.94
    "E31" = -1.46
    "F31" = -1.46
    "G31" = -0.16
    "H31" = 0.03999999999999999
    "I31" = -0.16
    "J31" = -1.26
    "K31" = -0.16
    "B32" = -0.35
    "C32" = 1.95
    "D32" = -1.45
    "E32" = -1.45
    "F32" = -0.15
    "G32" = 0.05
    "H32" = -0.15
    "I32" = -1.25
    "J32" = -0.15
    "K32" = 0.05
    "B33" = 1.96
    "C33" = -1.44
    "D33" = -1.44
    "E33" = -0.14
    "F33" = 0.05999999999999998
    "G33" = -0.14
    "H33" = -1.24
    "I33" = -0.14
    "J33" = 0.05999999999999998
    "K33" = -0.34
    "B34" = -1.67
    "C34" = -1.67
    "D34" = -0.37
    "E34" = -0.17
    "F34" = -0.37
    "G34" = -1.47
    "H34" = -0.37
    "I34" = -0.17
    "J34" = -0.5700000000000001
    "K34" = -0.27
    "B35" = -1.55
    "C35" = -0.25
    "D35" = -0.04999999999999999
    "E35" = -0.25
    "F35" = -1.35
    "G35" = -0.25
    "H35" = -0.04999999999999999
    "I35" = -0.45
    "J35" = -0.15
    "K35" = -0.55
    "B36" = -0.23
    "C36" = -0.03
    "D36" = -0.23
    "E36" = -1.33
    "F36" = -0.23
    "G36" = -0.03
    "H36" = -0.43
    "I36" = -0.13
    "J36" = -0.53
    "K36" = 0.27
    "B37" = 0.03999999999999999
    "C37" = -0.16
    "D37" = -1.26
    "E37" = -0.16
    "F37" = 0.03999999999999999
    "G37" = -0.36
    "H37" = -0.06000000000000002
    "I37" = -0.46
    "J37" = 0.34
    "K37" = 0.03999999999999999
    "B38" = -0.13
    "C38" = -1.23
    "D38" = -0.13
    "E38" = 0.07000000000000002
    "F38" = -0.33
    "G38" = -0.02999999999999998
    "H38" = -0.43
    "I38" = 0.3700000000000001
    "J38" = 0.07000000000000002
    "K38" = -0.02999999999999998
    "B39" = -1.24
    "C39" = -0.14
    "D39" = 0.06
    "E39" = -0.34
    "F39" = -0.04
    "G39" = -0.44
    "H39" = 0.36
    "I39" = 0.06
    "J39" = -0.04
    "K39" = -0.44
    "B40" = 0.07000000000000001
    "C40" = 0.27
    "D40" = -0.13
    "E40" = 0.17
    "F40" = -0.23
    "G40" = 0.5700000000000001
    "H40" = 0.27
    "I40" = 0.17
    "J40" = -0.23
    "K40" = -0.63
    "B41" = 0.14
    "C41" = -0.26
    "D41" = 0.04000000000000004
    "E41" = -0.36
    "F41" = 0.4400000000000001
    "G41" = 0.14
    "H41" = 0.04000000000000004
    "I41" = -0.36
    "J41" = -0.76
    "K41" = -0.76
    "B42" = -0.26
    "C42" = 0.04000000000000001
    "D42" = -0.36
    "E42" = 0.4400000000000001
    "F42" = 0.14
    "G42" = 0.04000000000000001
    "H42" = -0.36
    "I42" = -0.76
    "J42" = -0.76
    "K42" = 0.9400000000000001
    "B43" = 0.1
    "C43" = -0.3
    "D43" = 0.5
    "E43" = 0.2
    "F43" = 0.1
    "G43" = -0.3
    "H43" = -0.7000000000000001
    "I43" = -0.7000000000000001
    "J43" = 1
    "K43" = -0.4
    "B44" = -0.09000000000000002
    "C44" = 0.71
    "D44" = 0.41
    "E44" = 0.31
    "F44" = -0.09000000000000002
    "G44" = -0.49
    "H44" = -0.49
    "I44" = 1.21
    "J44" = -0.19
    "B45" = 0.6200000000000001
    "C45" = 0.3200000000000001
    "D45" = 0.2200000000000001
    "E45" = -0.1799999999999999
    "F45" = -0.58
    "G45" = -0.58
    "H45" = 1.12
    "I45" = -0.2799999999999999
    "B46" = 0.11
    "C46" = 0.009999999999999986
    "D46" = -0.39
    "E46" = -0.79
    "F46" = -0.79
    "G46" = 0.91
    "H46" = -0.49
    "B47" = 0.06000000000000001
    "C47" = -0.34
    "D47" = -0.74
    "E47" = -0.74
    "F47" = 0.9600000000000001
    "G47" = -0.44
    "B48" = -0.4
    "C48" = -0.8
    "D48" = -0.8
    "E48" = 0.9
    "F48" = -0.5
    "B49" = -0.79
    "C49" = -0.79
    "D49" = 0.91
    "E49" = -0.49
    "B50" = -0.67
    "C50" = 1.03
    "D50" = -0.37
    "B51" = 1
    "C51" = -0.4
    "B52" = -0.53
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
